# Update the "Price" (D) and "Volume(1h)" (E) columns on the active sheet
# to reflect a refreshed snapshot of crypto symbol data (GitHub Actions
# scheduled update). Values are written with a leading apostrophe so Excel
# stores them as literal text (matching the original "inlineStr" cells)
# instead of silently re-typing numeric-looking strings (e.g. "273.54",
# "1.26%") as numbers/percentages, which would both change the stored
# type and introduce floating point rounding.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'273.54"
$ws.Range("E2").Value = "'1.26%"
$ws.Range("E3").Value = "'0.23%"
$ws.Range("D4").Value = "'4.909"
$ws.Range("E4").Value = "'4.16%"
$ws.Range("D5").Value = "'0.06314"
$ws.Range("E5").Value = "'3.19%"
$ws.Range("D6").Value = "'6.919"
$ws.Range("E6").Value = "'2.66%"
$ws.Range("D7").Value = "'3.354"
$ws.Range("E7").Value = "'5.68%"
$ws.Range("D8").Value = "'1.323"
$ws.Range("E8").Value = "'47.38%"
$ws.Range("D9").Value = "'0.8829"
$ws.Range("E9").Value = "'3.09%"
$ws.Range("D10").Value = "'0.1473"
$ws.Range("E10").Value = "'2.85%"
$ws.Range("D11").Value = "'0.05093"
$ws.Range("E11").Value = "'2.29%"
$ws.Range("D12").Value = "'0.07385"
$ws.Range("E12").Value = "'3.96%"
$ws.Range("D13").Value = "'0.03180"
$ws.Range("E13").Value = "'0.18%"
$ws.Range("D14").Value = "'0.09047"
$ws.Range("E14").Value = "'0.18%"
$ws.Range("D15").Value = "'0.001563"
$ws.Range("E15").Value = "'2.11%"
$ws.Range("D16").Value = "'0.0006336"
$ws.Range("E16").Value = "'4.19%"
$ws.Range("D17").Value = "'0.006018"
$ws.Range("E17").Value = "'0.51%"
$ws.Range("D18").Value = "'3.473"
$ws.Range("E18").Value = "'0.31%"
$ws.Range("D21").Value = "'0.1336"
$ws.Range("E21").Value = "'4.37%"
$ws.Range("D22").Value = "'3.905"
$ws.Range("E22").Value = "'1.41%"
$ws.Range("D23").Value = "'0.04345"
$ws.Range("E23").Value = "'2.33%"
$ws.Range("D24").Value = "'0.001178"
$ws.Range("E24").Value = "'0.32%"
$ws.Range("D25").Value = "'0.003640"
$ws.Range("E25").Value = "'-12.37%"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'0.08%"
$ws.Range("D27").Value = "'0.0001702"
$ws.Range("E27").Value = "'1.29%"
$ws.Range("D40").Value = "'0.04051"
$ws.Range("E40").Value = "'2.58%"
$ws.Range("D41").Value = "'0.006596"
$ws.Range("E41").Value = "'57.55%"
$ws.Range("D42").Value = "'0.1161"
$ws.Range("E42").Value = "'3.64%"
$ws.Range("D43").Value = "'0.002221"
$ws.Range("E43").Value = "'9.10%"
$ws.Range("D44").Value = "'0.01262"
$ws.Range("E44").Value = "'-5.52%"
$ws.Range("D45").Value = "'0.00005343"
$ws.Range("E45").Value = "'4.08%"
$ws.Range("E46").Value = "'153.58%"
$ws.Range("D47").Value = "'0.02125"
$ws.Range("E47").Value = "'-13.17%"
